# Insert a new row for "Lesotho" (africa) above the existing "Liberia" row,
# keeping the list sorted alphabetically (row 114, pushing subsequent rows
# down by one).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(114).Insert()
$ws.Cells.Item(114, 1).Value = "Lesotho"
$ws.Cells.Item(114, 2).Value = "africa"

# Match the saved selection/scroll state from the authored workbook.
$ws.Range("A114").Select()
$excel.ActiveWindow.ScrollRow = 98
$excel.ActiveWindow.ScrollColumn = 1
